$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update file-name text values (WHO -> .csv, Domestic surveillance -> .csv.bz2) ---
# Write F13 first so the new shared-string entries land in the same order as the
# target workbook (Domestic surveillance ... .csv.bz2 before Canada ... .csv).
$ws.Range("F13").Value = "Domestic surveillance data - {suffix}_DISCOVER.csv.bz2"
$ws.Range("F3").Value = "Canada_COVID19_WHO_linelist-{suffix}_DISCOVER.csv"
$ws.Range("F4").Value = "Canada_COVID19_WHO_linelist-{suffix}_DISCOVER.csv"

# --- Highlight L3 / L4 with a solid yellow fill ---
$ws.Range("L3").Interior.Color = 65535
$ws.Range("L4").Interior.Color = 65535

# --- Update the sheet view / selection state ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F5").Select()
